$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.17578125
$ws.Range("C2").Value = 0.5859375
$ws.Range("J2").Value = 0.0078125
$ws.Range("P2").Value = 0.1328125
$ws.Range("S2").Value = 0.09765625
$ws.Range("C3").Value = 0.01986754966887417
$ws.Range("J3").Value = 0.01986754966887417
$ws.Range("P3").Value = 0.7152317880794702
$ws.Range("S3").Value = 0.2450331125827815
$ws.Range("P4").Value = 0.7666666666666667
$ws.Range("S4").Value = 0.2333333333333333
$ws.Range("B6").Value = 0.09883720930232558
$ws.Range("D6").Value = 0.005813953488372093
$ws.Range("F6").Value = 0.06395348837209303
$ws.Range("J6").Value = 0.2616279069767442
$ws.Range("O6").Value = 0.03488372093023256
$ws.Range("Q6").Value = 0.1569767441860465
$ws.Range("R6").Value = 0.06395348837209303
$ws.Range("S6").Value = 0.313953488372093
$ws.Range("B7").Value = 0.1666666666666667
$ws.Range("D7").Value = 0.02150537634408602
$ws.Range("E7").Value = 0.005376344086021506
$ws.Range("F7").Value = 0.05376344086021505
$ws.Range("J7").Value = 0.1182795698924731
$ws.Range("O7").Value = 0.005376344086021506
$ws.Range("Q7").Value = 0.1505376344086022
$ws.Range("R7").Value = 0.06451612903225806
$ws.Range("S7").Value = 0.4139784946236559
$ws.Range("B8").Value = 0.1133333333333333
$ws.Range("D8").Value = 0.01555555555555556
$ws.Range("F8").Value = 0.05333333333333334
$ws.Range("J8").Value = 0.09555555555555556
$ws.Range("O8").Value = 0.01111111111111111
$ws.Range("Q8").Value = 0.1911111111111111
$ws.Range("R8").Value = 0.09111111111111111
$ws.Range("S8").Value = 0.4288888888888889
$ws.Range("B9").Value = 0.09223300970873786
$ws.Range("D9").Value = 0.01456310679611651
$ws.Range("F9").Value = 0.03883495145631068
$ws.Range("J9").Value = 0.1019417475728155
$ws.Range("O9").Value = 0.01941747572815534
$ws.Range("Q9").Value = 0.1359223300970874
$ws.Range("R9").Value = 0.07766990291262135
$ws.Range("S9").Value = 0.5194174757281553
$ws.Range("B10").Value = 0.0947265625
$ws.Range("D10").Value = 0.0166015625
$ws.Range("E10").Value = 0.0009765625
$ws.Range("F10").Value = 0.068359375
$ws.Range("J10").Value = 0.1123046875
$ws.Range("O10").Value = 0.017578125
$ws.Range("Q10").Value = 0.1904296875
$ws.Range("R10").Value = 0.109375
$ws.Range("S10").Value = 0.3896484375
$ws.Range("G11").Value = 0.1136363636363636
$ws.Range("J11").Value = 0.08636363636363636
$ws.Range("K11").Value = 0.1409090909090909
$ws.Range("L11").Value = 0.6409090909090909
$ws.Range("S11").Value = 0.01818181818181818
$ws.Range("G12").Value = 0.8034682080924855
$ws.Range("J12").Value = 0.1098265895953757
$ws.Range("L12").Value = 0.06358381502890173
$ws.Range("S12").Value = 0.02312138728323699
$ws.Range("G13").Value = 0.64
$ws.Range("J13").Value = 0.24
$ws.Range("S13").Value = 0.12
$ws.Range("F15").Value = 0.01176470588235294
$ws.Range("H15").Value = 0.2411764705882353
$ws.Range("I15").Value = 0.08235294117647059
$ws.Range("J15").Value = 0.3
$ws.Range("K15").Value = 0.03529411764705882
$ws.Range("M15").Value = 0.01764705882352941
$ws.Range("O15").Value = 0.05882352941176471
$ws.Range("S15").Value = 0.2529411764705882
$ws.Range("F16").Value = 0.01829268292682927
$ws.Range("H16").Value = 0.1707317073170732
$ws.Range("I16").Value = 0.06707317073170732
$ws.Range("J16").Value = 0.4207317073170732
$ws.Range("K16").Value = 0.08536585365853659
$ws.Range("M16").Value = 0.01829268292682927
$ws.Range("O16").Value = 0.06097560975609756
$ws.Range("S16").Value = 0.1585365853658537
$ws.Range("F17").Value = 0.01657458563535912
$ws.Range("H17").Value = 0.2071823204419889
$ws.Range("I17").Value = 0.1270718232044199
$ws.Range("J17").Value = 0.3646408839779006
$ws.Range("K17").Value = 0.08839779005524862
$ws.Range("M17").Value = 0.03038674033149171
$ws.Range("N17").Value = 0.002762430939226519
$ws.Range("O17").Value = 0.05524861878453038
$ws.Range("S17").Value = 0.1077348066298343
$ws.Range("F18").Value = 0.02094240837696335
$ws.Range("H18").Value = 0.2146596858638743
$ws.Range("I18").Value = 0.1047120418848168
$ws.Range("J18").Value = 0.4083769633507853
$ws.Range("K18").Value = 0.06282722513089005
$ws.Range("M18").Value = 0.01047120418848168
$ws.Range("O18").Value = 0.05235602094240838
$ws.Range("S18").Value = 0.1256544502617801
$ws.Range("F19").Value = 0.01382886776145203
$ws.Range("H19").Value = 0.2221261884183232
$ws.Range("I19").Value = 0.09853068280034573
$ws.Range("J19").Value = 0.3569576490924806
$ws.Range("K19").Value = 0.1037165082108902
$ws.Range("M19").Value = 0.02592912705272256
$ws.Range("O19").Value = 0.05358686257562662
$ws.Range("S19").Value = 0.125324114088159
